# Insert 3 new data rows at row 255 (shifts existing rows 255-325 down to 258-328)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A255:A257").EntireRow.Insert()

# --- New row 255 ---
$ws.Cells.Item(255, 1).Value2  = 3
$ws.Cells.Item(255, 2).Value   = "Femacal de La Calera"
$ws.Cells.Item(255, 3).Value   = "Coquimbo"
$ws.Cells.Item(255, 4).Value2  = 44736
$ws.Cells.Item(255, 5).Value2  = 5
$ws.Cells.Item(255, 6).Value2  = 100112013
$ws.Cells.Item(255, 7).Value   = "Alcachofa"
$ws.Cells.Item(255, 8).Value   = "Argentina(o)"
$ws.Cells.Item(255, 9).Value   = "Primera"
$ws.Cells.Item(255, 10).Value2 = 70
$ws.Cells.Item(255, 11).Value2 = 18000
$ws.Cells.Item(255, 12).Value2 = 18000
$ws.Cells.Item(255, 13).Value2 = 18000
$ws.Cells.Item(255, 14).Value  = "$/caja 50 unidades"
$ws.Cells.Item(255, 15).Value  = "Provincia de Limarí"
$ws.Cells.Item(255, 16).Value2 = 360
$ws.Cells.Item(255, 17).Value2 = 50
$ws.Cells.Item(255, 18).Value  = "Hortaliza"

# --- New row 256 ---
$ws.Cells.Item(256, 1).Value2  = 3
$ws.Cells.Item(256, 2).Value   = "Femacal de La Calera"
$ws.Cells.Item(256, 3).Value   = "Coquimbo"
$ws.Cells.Item(256, 4).Value2  = 44736
$ws.Cells.Item(256, 5).Value2  = 5
$ws.Cells.Item(256, 6).Value2  = 100112013
$ws.Cells.Item(256, 7).Value   = "Alcachofa"
$ws.Cells.Item(256, 8).Value   = "Argentina(o)"
$ws.Cells.Item(256, 9).Value   = "Segunda"
$ws.Cells.Item(256, 10).Value2 = 75
$ws.Cells.Item(256, 11).Value2 = 17500
$ws.Cells.Item(256, 12).Value2 = 17500
$ws.Cells.Item(256, 13).Value2 = 17500
$ws.Cells.Item(256, 14).Value  = "$/caja 50 unidades"
$ws.Cells.Item(256, 15).Value  = "Provincia de Limarí"
$ws.Cells.Item(256, 16).Value2 = 350
$ws.Cells.Item(256, 17).Value2 = 50
$ws.Cells.Item(256, 18).Value  = "Hortaliza"

# --- New row 257 ---
$ws.Cells.Item(257, 1).Value2  = 3
$ws.Cells.Item(257, 2).Value   = "Femacal de La Calera"
$ws.Cells.Item(257, 3).Value   = "Coquimbo"
$ws.Cells.Item(257, 4).Value2  = 44736
$ws.Cells.Item(257, 5).Value2  = 5
$ws.Cells.Item(257, 6).Value2  = 100112013
$ws.Cells.Item(257, 7).Value   = "Alcachofa"
$ws.Cells.Item(257, 8).Value   = "Española"
$ws.Cells.Item(257, 9).Value   = "Extra"
$ws.Cells.Item(257, 10).Value2 = 95
$ws.Cells.Item(257, 11).Value2 = 21000
$ws.Cells.Item(257, 12).Value2 = 22000
$ws.Cells.Item(257, 13).Value2 = 21474
$ws.Cells.Item(257, 14).Value  = "$/caja 30 unidades"
$ws.Cells.Item(257, 15).Value  = "Provincia de Limarí"
$ws.Cells.Item(257, 16).Value2 = 716
$ws.Cells.Item(257, 17).Value2 = 30
$ws.Cells.Item(257, 18).Value  = "Hortaliza"
